$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the formatting of the last existing data row (33) onto the new row (34)
$ws.Range("A33:E33").Copy() | Out-Null
$ws.Range("A34:E34").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Cells.Item(34, 1).Value = 45938
$ws.Cells.Item(34, 2).Value = "21,5642"
$ws.Cells.Item(34, 3).Value = "15,3878"
$ws.Cells.Item(34, 4).Value = "15,3878"
$ws.Cells.Item(34, 5).Value = "15,3878"
